# Update electricity and stainless steel CF for #39
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GWP")

# Row 2: StainlessSteel - expected value updated, low/high become formulas
# derived from the expected (C) value instead of hard-coded numbers.
$ws.Range("C2").Value = 5.0231000000000003
$ws.Range("D2").Formula = '=$C2*0.75'
$ws.Range("E2").Formula = '=$C2*1.25'

# Row 12: Electricity - expected value updated, low/high become formulas
# derived from the expected (C) value instead of hard-coded numbers.
$ws.Range("C12").Value = 0.69711999999999996
$ws.Range("D12").Formula = '=$C12*0.75'
$ws.Range("E12").Formula = '=$C12*1.25'

# Reflect the author leaving the workbook with the GWP sheet active and
# cell G14 selected (the "info" sheet was previously the active/selected one).
$ws.Activate()
$ws.Range("G14").Select()
